# Corrected tests to precreate containers into the container handler hash
# to avoid attempts to access via solr.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("testTopLinkerUpload")

# Top Container Record No. values bumped by one (pre-created container ids)
$ws.Range("M6").Value = 54556
$ws.Range("M7").Value = 54556
$ws.Range("M8").Value = 54556

# Location Record No. values now point at a pre-created location (id 4)
# instead of the old solr-lookup ids
$ws.Range("J9").Value = 4
$ws.Range("J10").Value = 4

# Update the active selection to match where the author left off editing
$ws.Range("M7:M8").Select()

$wb.Save()
